# Fruta / hortaliza, semanal
# Insert two new daily price records for "Cebolla" right before the current
# row 1061 (pushing the existing 1061-1114 block down to 1063-1116), then
# populate the two freshly-inserted rows with the new observations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 1061, shifting everything at/after 1061 down by 2.
$ws.Range("A1061:R1062").EntireRow.Insert()

# New row 1061: "1a nueva(o)" onion lot, O'Higgins region.
$ws.Range("A1061").Value = 5
$ws.Range("B1061").Value = "Macroferia Regional de Talca"
$ws.Range("C1061").Value = "Maule"
$ws.Range("D1061").Value = 45267
$ws.Range("E1061").Value = 7
$ws.Range("F1061").Value = 100112004
$ws.Range("G1061").Value = "Cebolla"
$ws.Range("H1061").Value = "Sin especificar"
$ws.Range("I1061").Value = "1a nueva(o)"
$ws.Range("J1061").Value = 50000
$ws.Range("K1061").Value = 2800
$ws.Range("L1061").Value = 2800
$ws.Range("M1061").Value = 2800
$ws.Range("N1061").Value = "$/paquete 10 unidades (volumen en unidades)"
$ws.Range("O1061").Value = "Región de O'Higgins"
$ws.Range("P1061").Value = 280
$ws.Range("Q1061").Value = 10
$ws.Range("R1061").Value = "Hortaliza"

# New row 1062: "2a nueva(o)" onion lot, O'Higgins region.
$ws.Range("A1062").Value = 5
$ws.Range("B1062").Value = "Macroferia Regional de Talca"
$ws.Range("C1062").Value = "Maule"
$ws.Range("D1062").Value = 45267
$ws.Range("E1062").Value = 7
$ws.Range("F1062").Value = 100112004
$ws.Range("G1062").Value = "Cebolla"
$ws.Range("H1062").Value = "Sin especificar"
$ws.Range("I1062").Value = "2a nueva(o)"
$ws.Range("J1062").Value = 15000
$ws.Range("K1062").Value = 1800
$ws.Range("L1062").Value = 1800
$ws.Range("M1062").Value = 1800
$ws.Range("N1062").Value = "$/paquete 10 unidades (volumen en unidades)"
$ws.Range("O1062").Value = "Región de O'Higgins"
$ws.Range("P1062").Value = 180
$ws.Range("Q1062").Value = 10
$ws.Range("R1062").Value = "Hortaliza"
